$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    3  = 10.95594069999788
    4  = 11.30781500000012
    5  = 11.31903919999604
    6  = 11.44202339999902
    7  = 11.7871535000013
    8  = 10.84841809999489
    9  = 10.06526139999914
    10 = 10.51989439999306
    11 = 9.429674299994076
    12 = 10.11985519999871
    13 = 10.13611189999938
    14 = 9.350660700001754
    15 = 10.32756700000027
    16 = 9.947236100000737
    17 = 9.556730099997367
    18 = 25.48973509999632
    19 = 5.525759299998754
    20 = 5.760017399996286
    21 = 5.209805399994366
    22 = 5.331813099997817
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
